# Updated cryptos list on Mon Dec 25 04:37:47 UTC 2023 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns of the crypto table, and
# swaps the NEARProtocol/VeChain rows (38-39) back to their new ranking order.
#
# Price cells are prefixed with a leading apostrophe so Excel stores them as
# text (matching the sheet's original inline-string cells) instead of
# silently parsing e.g. "0.608" as a number; .Style is reset to "Normal"
# right after so no stray number-format style is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.139.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "'2.279.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'112.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.31%  "
$ws.Range("D6").Value = "'265.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("D10").Value = "'47.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "'8.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "'15.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").Value = "'2.622.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "'2.277.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "'43.245.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").Value = "'71.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").Value = "'2.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").Value = "'232.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").Value = "'9.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("D27").Value = "'11.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").Value = "'40.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.16%  "
$ws.Range("D30").Value = "'3.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.12%  "
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").Value = "'172.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.34%  "
$ws.Range("D33").Value = "'21.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.67%  "
$ws.Range("D34").Value = "'0.0907"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.37%  "
$ws.Range("D35").Value = "'5.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.74%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Value = "'4.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.28%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'3.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0355"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  -5.69%  "
$ws.Range("E41").Value = "  +11.07%  "
$ws.Range("D42").Value = "'76.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.73%  "
$ws.Range("D43").Value = "'13.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.36%  "
$ws.Range("E44").Value = "  -5.40%  "
$ws.Range("D45").Value = "'6.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.51%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").Value = "'8.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").Value = "'103.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("D50").Value = "'0.0994"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("E51").Value = "  +0.99%  "
